$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Antenna named value (Sheet1!$T$13) used by the effective-range formulas.
$ws.Range("T13").Value = 25000000000

# New "Radius" / "Flyby" reference figures added near the bottom-right of the sheet.
$ws.Range("L21").Formula = "=1550800"
$ws.Range("M21").Value = "Radius"

$ws.Range("L22").Value = 9700000
$ws.Range("M22").Value = "Flyby"

$ws.Range("L23").Formula = "=L21*4"
$ws.Range("L24").Formula = "=L23/1000"

# Leave the selection where the author left it when saving.
$ws.Range("L25").Select()
